$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks numeric would otherwise be auto-coerced to a
# number by Excel's smart-typing. For those, force Text format first, write
# the value, then reset the style back to Normal/default afterwards so no
# lingering number-format/style change is left on the cell.

$ws.Range('D2').Value = '27.526.84'
$ws.Range('E2').Value = '  +2.10%  '
$ws.Range('D3').Value = '1.564.92'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  -1.52%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '210.27'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.488'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.988'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.52%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.42'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0868'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.52%  '
$ws.Range('D12').Value = '1.791.99'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = '1.606.26'
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.75'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.520'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('D16').Value = '27.517.54'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '62.01'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '224.33'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +4.04%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.52'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('D20').Value = '0.0₃0704'
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('E21').Value = '  -1.43%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.14'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.37'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.94'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '150.04'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.62%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '15.19'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.21%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.62'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.107'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.13'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0470'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.23'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('D33').Value = '1.464.38'
$ws.Range('E33').Value = '  +3.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.18'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.10'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.62'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.64%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.31'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.540'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.79%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.815'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.94'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +11.87%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.72'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.52%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.34'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.04%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.988'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.973'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.51%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '65.11'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('D47').Value = '1.705.53'
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '86.30'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('D50').Value = '0.0₆0101'
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0952'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.59%  '
